$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1764
$ws.Range("J17").Value = 1999.3846
$ws.Range("L17").Value = 5998.1538
$ws.Range("N17").Value = -6334.1538
$ws.Range("H28").Value = 791
$ws.Range("I28").Value = 863.75
$ws.Range("K28").Value = 863.75
$ws.Range("M28").Value = -378.75
$ws.Range("H40").Value = 2500
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4350
$ws.Range("H41").Value = 396
$ws.Range("I41").Value = 383.6
$ws.Range("K41").Value = 383.6
$ws.Range("M41").Value = 56.39999999999998
$ws.Range("H53").Value = 331.66666
$ws.Range("I53").Value = 483.33334
$ws.Range("J53").Value = 255.83333
$ws.Range("K53").Value = 483.33334
$ws.Range("L53").Value = 255.83333
$ws.Range("M53").Value = 153.66666
$ws.Range("N53").Value = -1529.83333
$ws.Range("H70").Value = 60390
$ws.Range("J70").Value = 117066
$ws.Range("L70").Value = 351198
$ws.Range("N70").Value = -351738
$ws.Range("H73").Value = 60390
$ws.Range("J73").Value = 117066
$ws.Range("L73").Value = 351198
$ws.Range("N73").Value = -353070
$ws.Range("H76").Value = 6425.4736
$ws.Range("I76").Value = 5509.8887
$ws.Range("K76").Value = 5509.8887
$ws.Range("M76").Value = -5194.8887
$ws.Range("H79").Value = 6425.4736
$ws.Range("I79").Value = 5509.8887
$ws.Range("K79").Value = 5509.8887
$ws.Range("M79").Value = -4417.8887
$ws.Range("H86").Value = 4443.3887
$ws.Range("I86").Value = 3635.2222
$ws.Range("J86").Value = 5251.5557
$ws.Range("K86").Value = 3635.2222
$ws.Range("L86").Value = 5251.5557
$ws.Range("M86").Value = -2512.2222
$ws.Range("N86").Value = -7497.5557
$ws.Range("H89").Value = 4443.3887
$ws.Range("I89").Value = 3635.2222
$ws.Range("J89").Value = 5251.5557
$ws.Range("K89").Value = 18176.111
$ws.Range("L89").Value = 26257.7785
$ws.Range("M89").Value = -12560.111
$ws.Range("N89").Value = -37489.7785
$ws.Range("H92").Value = 1530.1666
$ws.Range("J92").Value = 1971.75
$ws.Range("L92").Value = 1971.75
$ws.Range("N92").Value = -4467.75
$ws.Range("H98").Value = 1728.15
$ws.Range("I98").Value = 1661.2106
$ws.Range("K98").Value = 1661.2106
$ws.Range("M98").Value = -163.2106000000001
$ws.Range("H106").Value = 22855.416
$ws.Range("I106").Value = 25626.5
$ws.Range("K106").Value = 25626.5
$ws.Range("M106").Value = -24995.5
$ws.Range("H107").Value = 547.36365
$ws.Range("I107").Value = 527.55554
$ws.Range("J107").Value = 636.5
$ws.Range("K107").Value = 527.55554
$ws.Range("L107").Value = 636.5
$ws.Range("M107").Value = 1392.44446
$ws.Range("N107").Value = -4476.5
$ws.Range("H122").Value = 1728.15
$ws.Range("I122").Value = 1661.2106
$ws.Range("K122").Value = 4983.6318
$ws.Range("M122").Value = -2533.6318
$ws.Range("H137").Value = 2081.8965
$ws.Range("I137").Value = 1032.6666
$ws.Range("J137").Value = 3798.818
$ws.Range("K137").Value = 3097.9998
$ws.Range("L137").Value = 11396.454
$ws.Range("M137").Value = -547.9998
$ws.Range("N137").Value = -16496.454

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1197.5
$ws.Range("I2").Value = 1056.2
$ws.Range("J2").Value = 1433
$ws.Range("K2").Value = 1056.2
$ws.Range("L2").Value = 1433
$ws.Range("M2").Value = -943.2
$ws.Range("N2").Value = -1659
$ws.Range("H32").Value = 15209.926
$ws.Range("J32").Value = 23399.824
$ws.Range("L32").Value = 23399.824
$ws.Range("N32").Value = -23973.824
$ws.Range("H116").Value = 1197.5
$ws.Range("I116").Value = 1056.2
$ws.Range("J116").Value = 1433
$ws.Range("K116").Value = 1056.2
$ws.Range("L116").Value = 1433
$ws.Range("M116").Value = 1237.8
$ws.Range("N116").Value = -6021

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1197.5
$ws.Range("I3").Value = 1056.2
$ws.Range("J3").Value = 1433
$ws.Range("K3").Value = 1056.2
$ws.Range("L3").Value = 1433
$ws.Range("M3").Value = -942.2
$ws.Range("N3").Value = -1661
$ws.Range("H7").Value = 25010000
$ws.Range("I7").Value = 25010000
$ws.Range("K7").Value = 25010000
$ws.Range("M7").Value = -25009887
$ws.Range("H86").Value = 3903.5
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 6807
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 6807
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -9053
$ws.Range("H89").Value = 3903.5
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 6807
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 34035
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -45267
$ws.Range("H105").Value = 4542.6294
$ws.Range("I105").Value = 4003.7144
$ws.Range("K105").Value = 4003.7144
$ws.Range("M105").Value = -2256.7144

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1625.25
$ws.Range("I2").Value = 1042.4286
$ws.Range("J2").Value = 5705
$ws.Range("K2").Value = 1042.4286
$ws.Range("L2").Value = 5705
$ws.Range("M2").Value = -929.4286
$ws.Range("N2").Value = -5931
$ws.Range("H17").Value = 1007
$ws.Range("I17").Value = 1007
$ws.Range("K17").Value = 1007
$ws.Range("M17").Value = -833
$ws.Range("H33").Value = 1988.6
$ws.Range("I33").Value = 985.75
$ws.Range("J33").Value = 6000
$ws.Range("K33").Value = 985.75
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = -606.75
$ws.Range("N33").Value = -6758
$ws.Range("H42").Value = 3032.6667
$ws.Range("I42").Value = 49
$ws.Range("K42").Value = 49
$ws.Range("M42").Value = 544
$ws.Range("H44").Value = 11021.667
$ws.Range("I44").Value = 65
$ws.Range("J44").Value = 16500
$ws.Range("K44").Value = 65
$ws.Range("L44").Value = 16500
$ws.Range("M44").Value = 377
$ws.Range("N44").Value = -17384
$ws.Range("H55").Value = 75
$ws.Range("I55").Value = 75
$ws.Range("K55").Value = 75
$ws.Range("M55").Value = 240
$ws.Range("H99").Value = 11722.305
$ws.Range("I99").Value = 7708.467
$ws.Range("K99").Value = 7708.467
$ws.Range("M99").Value = -6210.467
$ws.Range("H122").Value = 9102.4
$ws.Range("I122").Value = 9166
$ws.Range("J122").Value = 9007
$ws.Range("K122").Value = 27498
$ws.Range("L122").Value = 27021
$ws.Range("M122").Value = -25048
$ws.Range("N122").Value = -31921
$ws.Range("H126").Value = 11722.305
$ws.Range("I126").Value = 7708.467
$ws.Range("K126").Value = 23125.401
$ws.Range("M126").Value = -20655.401
$ws.Range("H132").Value = 3096.9583
$ws.Range("I132").Value = 2650.1875
$ws.Range("K132").Value = 7950.5625
$ws.Range("M132").Value = -5420.5625
$ws.Range("H134").Value = 2003.0286
$ws.Range("I134").Value = 1244.7142
$ws.Range("K134").Value = 3734.1426
$ws.Range("M134").Value = -1199.1426

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 73001260
$ws.Range("I4").Value = 127750350
$ws.Range("J4").Value = 2491
$ws.Range("K4").Value = 383251050
$ws.Range("L4").Value = 7473
$ws.Range("M4").Value = -383250938
$ws.Range("N4").Value = -7697
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H102").Value = 1157.3513
$ws.Range("I102").Value = 543.12
$ws.Range("K102").Value = 543.12
$ws.Range("M102").Value = 1078.88
$ws.Range("H122").Value = 502689.88
$ws.Range("I122").Value = 64199.812
$ws.Range("K122").Value = 192599.436
$ws.Range("M122").Value = -190149.436

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10050000
$ws.Range("J2").Value = 99999
$ws.Range("L2").Value = 99999
$ws.Range("N2").Value = -100223
$ws.Range("H40").Value = 1956.1666
$ws.Range("I40").Value = 1956.1666
$ws.Range("K40").Value = 1956.1666
$ws.Range("M40").Value = -1820.1666
$ws.Range("H132").Value = 5376.9443
$ws.Range("I132").Value = 4479.1
$ws.Range("K132").Value = 13437.3
$ws.Range("M132").Value = -10907.3

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 3620.5
$ws.Range("J41").Value = 3620.5
$ws.Range("L41").Value = 3620.5
$ws.Range("N41").Value = -4400.5
$ws.Range("H96").Value = 2082.75
$ws.Range("I96").Value = 1777
$ws.Range("K96").Value = 1777
$ws.Range("M96").Value = -404
$ws.Range("H122").Value = 1559.125
$ws.Range("I122").Value = 1424.7858
$ws.Range("K122").Value = 4274.357400000001
$ws.Range("M122").Value = -1824.357400000001
$ws.Range("H136").Value = 25675.146
$ws.Range("I136").Value = 1151.2693
$ws.Range("J136").Value = 68183.2
$ws.Range("K136").Value = 3453.8079
$ws.Range("L136").Value = 204549.6
$ws.Range("M136").Value = -903.8078999999998
$ws.Range("N136").Value = -209649.6
